# Add a new row (16) at the bottom of the error-count log:
#   Date = 11/19/2025 (serial 45980), Error Count = 3
# Mirrors the existing rows: copy the date cell's format from the row
# above (A15) so the new date cell picks up the same style (short-date
# number format) instead of minting a brand-new style entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 15
$newRow = $lastRow + 1

$ws.Range("A$lastRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A$newRow").Value = 45980
$ws.Range("B$newRow").Value = 3

$ws.Range("A$newRow`:B$newRow").Select() | Out-Null
